# Adding the changes we made on may 9th
#
# The accelerometer data table (x, y, z in columns A:B:C, header in row 1)
# gets 4 new sample rows inserted right after the header (pushing the
# existing 20 rows of data down by 4), and 6 new sample rows appended at
# the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows directly below the header (old row 2 becomes row 6) ---
$ws.Rows("2:5").Insert()
# The insert copies formatting from the row above (the bold header); clear
# it so the new data rows look like ordinary, unstyled data rows.
$ws.Rows("2:5").ClearFormats()

$newTopRows = @(
    @(-0.4100122451782226, 0.3157248497009277, -0.4111190438270569),
    @(-0.3957743644714355, 0.3475203514099121, -0.2738307416439056),
    @(-0.5445261001586914, 0.2971320152282715, -0.1380582749843597),
    @(-0.0015645027160644, 0.2662014961242676, -0.2422323226928711)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value2 = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $newTopRows[$i][2]
}

# --- Append 6 new rows after the (now shifted) existing data, rows 26-31 ---
$newBottomRows = @(
    @(1.732457160949707, 0.4553084373474121, 2.453210830688477),
    @(-0.187225341796875, 0.2255609035491943, -0.6162976026535034),
    @(-0.9302024841308594, 0.2106423377990722, -0.1286094188690185),
    @(-1.170828819274902, -0.4288506507873535, -0.3933718204498291),
    @(-0.678126335144043, 0.426605224609375, -0.0270633697509765),
    @(-0.493565559387207, 0.1020381450653076, 0.0103309154510498)
)

for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = 26 + $i
    $ws.Cells.Item($r, 1).Value2 = $newBottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $newBottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $newBottomRows[$i][2]
}
